$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated main: refreshed the generated project cost figures in row 2.
$ws.Range("A2").Value = 7073
$ws.Range("B2").Value = 6634
$ws.Range("C2").Value = 3955
$ws.Range("D2").Value = 5012
$ws.Range("E2").Value = 6380
$ws.Range("F2").Value = 3413
